$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ROI geometries to smaller, downloadable-size bounding boxes,
# and rename "hungry-romania" -> "hungry" / "bulgaria-romania" -> "bulgaria".
$ws.Range("A3").Value = "  [[[-0.09707574480274594, 52.29419689960738],`n           [-0.09707574480274594, 51.89598091800008],`n           [1.064728454416004, 51.89598091800008],`n           [1.064728454416004, 52.29419689960738]]]"
$ws.Range("A4").Value = "[[[1.2866262035756648, 48.385430020942366],`n          [1.2866262035756648, 47.93473449186214],`n          [2.445683820763165, 47.93473449186214],`n          [2.445683820763165, 48.385430020942366]]]"
$ws.Range("A5").Value = "[[[[20.33654768582722, 46.81969006391422],`n           [20.33654768582722, 46.36293233084845],`n           [21.50659163113972, 46.36293233084845],`n           [21.50659163113972, 46.81969006391422]]]"
$ws.Range("D5").Value = "hungry"
$ws.Range("A7").Value = " [[[-5.73566016354585, 42.21381966999416],`n           [-5.73566016354585, 41.70731558747302],`n           [-4.57110938229585, 41.70731558747302],`n           [-4.57110938229585, 42.21381966999416]]]"
$ws.Range("A8").Value = "[[[27.303049105663984, 44.042971343838296],`n           [27.303049105663984, 43.55928538069457],`n           [28.445627230663984, 43.55928538069457],`n           [28.445627230663984, 44.042971343838296]]]"
$ws.Range("D8").Value = "bulgaria"
$ws.Range("A9").Value = " [[[23.559466433828824, 56.32475801105907],`n           [23.559466433828824, 55.94212139772337],`n           [24.740496707266324, 55.94212139772337],`n           [24.740496707266324, 56.32475801105907]]]"

# Restore the active selection to C9, matching the saved workbook state.
$ws.Range("C9").Select() | Out-Null
